# Add new issue #20 as a new row (row 16) at the bottom of the Issues sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("Issue ID") holds numeric-looking text (e.g. "19") stored as text,
# not as a number, in the source data. Format the cell as text before writing
# the value so "20" is preserved as a string, then reset the style back to
# the default ("Normal") so no extra formatting is introduced.
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "20"
$ws.Range("A16").Style = "Normal"

$ws.Range("B16").Value = "FR_HIRING"
$ws.Range("C16").Value = "open"
$ws.Range("D16").Value = "2025-03-26T06:40:25Z"
$ws.Range("E16").Value = "bug"
